# Update computed "loading_percent" results for the 380 kV case (res_line sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 13.12240423021371
$ws.Range("D2").Value = 5.264534762822142
$ws.Range("E2").Value = 13.56397927401664
$ws.Range("F2").Value = 24.76176567724934
$ws.Range("G2").Value = 29.54827907061934
$ws.Range("H2").Value = 14.58574271590618
$ws.Range("K2").Value = 14.18055352004724
$ws.Range("L2").Value = 9.308426832388779
$ws.Range("M2").Value = 16.83976888856147
$ws.Range("O2").Value = 22.28874921790268

# Row 3
$ws.Range("C3").Value = 13.10207764057802
$ws.Range("D3").Value = 5.206238087509748
$ws.Range("E3").Value = 13.59955362221075
$ws.Range("F3").Value = 24.86674808705366
$ws.Range("G3").Value = 29.73053770777309
$ws.Range("H3").Value = 14.65685682476672
$ws.Range("K3").Value = 13.58454316405913
$ws.Range("L3").Value = 9.338282633960173
$ws.Range("M3").Value = 16.59111066583316
$ws.Range("O3").Value = 22.41597124582423

# Row 4
$ws.Range("C4").Value = 13.09281696559318
$ws.Range("D4").Value = 5.16983364272269
$ws.Range("E4").Value = 13.62417269842451
$ws.Range("F4").Value = 24.93933869123804
$ws.Range("G4").Value = 29.85466405121399
$ws.Range("H4").Value = 14.70342966797513
$ws.Range("K4").Value = 13.20461045418747
$ws.Range("L4").Value = 9.357743079053595
$ws.Range("M4").Value = 16.43801705684716
$ws.Range("O4").Value = 22.50011331733047

# Row 5
$ws.Range("C5").Value = 13.08985606462192
$ws.Range("D5").Value = 5.154854019836776
$ws.Range("E5").Value = 13.63490221062579
$ws.Range("F5").Value = 24.97095463385642
$ws.Range("G5").Value = 29.90829721462185
$ws.Range("H5").Value = 14.72313966011784
$ws.Range("K5").Value = 13.04645672469936
$ws.Range("L5").Value = 9.365957687355252
$ws.Range("M5").Value = 16.37559079866897
$ws.Range("O5").Value = 22.53591353893065

# Row 6
$ws.Range("C6").Value = 13.08941358649544
$ws.Range("D6").Value = 5.152358198618626
$ws.Range("E6").Value = 13.6367259041179
$ws.Range("F6").Value = 24.9763270298875
$ws.Range("G6").Value = 29.91738656952582
$ws.Range("H6").Value = 14.726456646018
$ws.Range("K6").Value = 13.0200003833113
$ws.Range("L6").Value = 9.367338902511959
$ws.Range("M6").Value = 16.36522451029147
$ws.Range("O6").Value = 22.54194932360014

# Row 7
$ws.Range("C7").Value = 13.09277373872127
$ws.Range("D7").Value = 5.169632195606972
$ws.Range("E7").Value = 13.62431457958695
$ws.Range("F7").Value = 24.93975684994214
$ws.Range("G7").Value = 29.85537504387678
$ws.Range("H7").Value = 14.70369252353282
$ws.Range("K7").Value = 13.20249074685065
$ws.Range("L7").Value = 9.357852712328627
$ws.Range("M7").Value = 16.43717522710828
$ws.Range("O7").Value = 22.50059001625291

# Row 8
$ws.Range("C8").Value = 13.11472946905769
$ws.Range("D8").Value = 5.244565970790762
$ws.Range("E8").Value = 13.57566846755983
$ws.Range("F8").Value = 24.7962701728902
$ws.Range("G8").Value = 29.60857278977065
$ws.Range("H8").Value = 14.60965908875011
$ws.Range("K8").Value = 13.97805540495842
$ws.Range("L8").Value = 9.318487092550287
$ws.Range("M8").Value = 16.75415694634613
$ws.Range("O8").Value = 22.33136203456446

# Row 9
$ws.Range("C9").Value = 13.18316066743232
$ws.Range("D9").Value = 5.386288943890692
$ws.Range("E9").Value = 13.50234751162623
$ws.Range("F9").Value = 24.5798356762701
$ws.Range("G9").Value = 29.22250643203517
$ws.Range("H9").Value = 14.44834564076516
$ws.Range("K9").Value = 15.38099790412431
$ws.Range("L9").Value = 9.250227462131463
$ws.Range("M9").Value = 17.36962752307917
$ws.Range("O9").Value = 22.04750780050535

# Row 10
$ws.Range("C10").Value = 13.24863063812148
$ws.Range("D10").Value = 5.486732737809844
$ws.Range("E10").Value = 13.46199303183114
$ws.Range("F10").Value = 24.46096937090096
$ws.Range("G10").Value = 28.99984242453252
$ws.Range("H10").Value = 14.34391382677887
$ws.Range("K10").Value = 16.33197132361361
$ws.Range("L10").Value = 9.205495919608421
$ws.Range("M10").Value = 17.81438435122779
$ws.Range("O10").Value = 21.86846855061031

# Row 11
$ws.Range("C11").Value = 13.28164011387624
$ws.Range("D11").Value = 5.531526927877612
$ws.Range("E11").Value = 13.44657961280424
$ws.Range("F11").Value = 24.41572169831113
$ws.Range("G11").Value = 28.912054612992
$ws.Range("H11").Value = 14.29946673488286
$ws.Range("K11").Value = 16.74598982011367
$ws.Range("L11").Value = 9.186316983583408
$ws.Range("M11").Value = 18.01435357322966
$ws.Range("O11").Value = 21.79348046902464

# Row 12
$ws.Range("C12").Value = 13.29459662219622
$ws.Range("D12").Value = 5.548351947941208
$ws.Range("E12").Value = 13.44116696433391
$ws.Range("F12").Value = 24.3998648722204
$ws.Range("G12").Value = 28.88077456689186
$ws.Range("H12").Value = 14.28307614880895
$ws.Range("K12").Value = 16.90000964944301
$ws.Range("L12").Value = 9.179222182472735
$ws.Range("M12").Value = 18.08967859826507
$ws.Range("O12").Value = 21.76601773602575

# Row 13
$ws.Range("C13").Value = 13.29178602796945
$ws.Range("D13").Value = 5.544734628166508
$ws.Range("E13").Value = 13.44231380365166
$ws.Range("F13").Value = 24.40322297878474
$ws.Range("G13").Value = 28.88742363040533
$ws.Range("H13").Value = 14.28658655041488
$ws.Range("K13").Value = 16.86696272442932
$ws.Range("L13").Value = 9.180742715345035
$ws.Range("M13").Value = 18.07347469027449
$ws.Range("O13").Value = 21.77189072333601

# Row 14
$ws.Range("C14").Value = 13.28269694917713
$ws.Range("D14").Value = 5.532913936658572
$ws.Range("E14").Value = 13.44612580802718
$ws.Range("F14").Value = 24.41439149871881
$ws.Range("G14").Value = 28.90944171636477
$ws.Range("H14").Value = 14.29810943830173
$ws.Range("K14").Value = 16.75871689153558
$ws.Range("L14").Value = 9.185729928126998
$ws.Range("M14").Value = 18.02055894306103
$ws.Range("O14").Value = 21.7912023533131

# Row 15
$ws.Range("C15").Value = 13.27718884699267
$ws.Range("D15").Value = 5.525655261520999
$ws.Range("E15").Value = 13.44851601799503
$ws.Range("F15").Value = 24.42139914645269
$ws.Range("G15").Value = 28.92318473687233
$ws.Range("H15").Value = 14.30522493307913
$ws.Range("K15").Value = 16.69205139634465
$ws.Range("L15").Value = 9.188806588160379
$ws.Range("M15").Value = 17.98809280346714
$ws.Range("O15").Value = 21.8031530195387

# Row 16
$ws.Range("C16").Value = 13.24653769968183
$ws.Range("D16").Value = 5.483786533271576
$ws.Range("E16").Value = 13.46305964164538
$ws.Range("F16").Value = 24.46410473214496
$ws.Range("G16").Value = 29.00585314002691
$ws.Range("H16").Value = 14.34688016477576
$ws.Range("K16").Value = 16.30453247162543
$ws.Range("L16").Value = 9.206772814669844
$ws.Range("M16").Value = 17.80126330715467
$ws.Range("O16").Value = 21.87349957157353

# Row 17
$ws.Range("C17").Value = 13.22855558075676
$ws.Range("D17").Value = 5.457865308251641
$ws.Range("E17").Value = 13.47273627870281
$ws.Range("F17").Value = 24.49256975080226
$ws.Range("G17").Value = 29.06004155150934
$ws.Range("H17").Value = 14.37321832849858
$ws.Range("K17").Value = 16.06197547941018
$ws.Range("L17").Value = 9.218093839981471
$ws.Range("M17").Value = 17.68600292246513
$ws.Range("O17").Value = 21.91831243609457

# Row 18
$ws.Range("C18").Value = 13.21851696010202
$ws.Range("D18").Value = 5.442872049216268
$ws.Range("E18").Value = 13.47857909472395
$ws.Range("F18").Value = 24.50977227728478
$ws.Range("G18").Value = 29.09247902341734
$ws.Range("H18").Value = 14.38865527564468
$ws.Range("K18").Value = 15.92071817200501
$ws.Range("L18").Value = 9.224715509324218
$ws.Range("M18").Value = 17.61949021174735
$ws.Range("O18").Value = 21.94469514320534

# Row 19
$ws.Range("C19").Value = 13.21517052236765
$ws.Range("D19").Value = 5.437781413209611
$ws.Range("E19").Value = 13.48060493476374
$ws.Range("F19").Value = 24.51573907879481
$ws.Range("G19").Value = 29.10367920427413
$ws.Range("H19").Value = 14.39393139517968
$ws.Range("K19").Value = 15.87259408065259
$ws.Range("L19").Value = 9.226976415990991
$ws.Range("M19").Value = 17.59693459873215
$ws.Range("O19").Value = 21.95373206730625

# Row 20
$ws.Range("C20").Value = 13.23043837460394
$ws.Range("D20").Value = 5.460633428137858
$ws.Range("E20").Value = 13.4716775027665
$ws.Range("F20").Value = 24.4894536154199
$ws.Range("G20").Value = 29.05414153197786
$ws.Range("H20").Value = 14.37038478177209
$ws.Range("K20").Value = 16.08797735899107
$ws.Range("L20").Value = 9.216877303586637
$ws.Range("M20").Value = 17.69829559488701
$ws.Range("O20").Value = 21.91347911670877

# Row 21
$ws.Range("C21").Value = 13.2853543046188
$ws.Range("D21").Value = 5.536389759497016
$ws.Range("E21").Value = 13.44499461533238
$ws.Range("F21").Value = 24.4110762987445
$ws.Range("G21").Value = 28.90292101910529
$ws.Range("H21").Value = 14.2947129249605
$ws.Range("K21").Value = 16.79058685209231
$ws.Range("L21").Value = 9.18426051030943
$ws.Range("M21").Value = 18.03611287990551
$ws.Range("O21").Value = 21.7855046809356

# Row 22
$ws.Range("C22").Value = 13.323902145942
$ws.Range("D22").Value = 5.585095347111268
$ws.Range("E22").Value = 13.43002791565628
$ws.Range("F22").Value = 24.36730185384194
$ws.Range("G22").Value = 28.81554275905979
$ws.Range("H22").Value = 14.24782524318105
$ws.Range("K22").Value = 17.23366414060615
$ws.Range("L22").Value = 9.163921727138092
$ws.Range("M22").Value = 18.25454409689032
$ws.Range("O22").Value = 21.70731023992661

# Row 23
$ws.Range("C23").Value = 13.3030879551779
$ws.Range("D23").Value = 5.559176673467823
$ws.Range("E23").Value = 13.43778951645375
$ws.Range("F23").Value = 24.38998081867231
$ws.Range("G23").Value = 28.86112327722689
$ws.Range("H23").Value = 14.27261487292147
$ws.Range("K23").Value = 16.99868578164676
$ws.Range("L23").Value = 9.17468752892777
$ws.Range("M23").Value = 18.13819783196013
$ws.Range("O23").Value = 21.7485442663782

# Row 24
$ws.Range("C24").Value = 13.2295862299231
$ws.Range("D24").Value = 5.459382243316125
$ws.Range("E24").Value = 13.47215530453035
$ws.Range("F24").Value = 24.49085981175996
$ws.Range("G24").Value = 29.05680493336271
$ws.Range("H24").Value = 14.37166491003233
$ws.Range("K24").Value = 16.07622753416843
$ws.Range("L24").Value = 9.217426947481981
$ws.Range("M24").Value = 17.69273884544696
$ws.Range("O24").Value = 21.91566233199712

# Row 25
$ws.Range("C25").Value = 13.16195720951999
$ws.Range("D25").Value = 5.348560206698385
$ws.Range("E25").Value = 13.51981326183449
$ws.Range("F25").Value = 24.63137507834882
$ws.Range("G25").Value = 29.31633069208215
$ws.Range("H25").Value = 14.48951266327271
$ws.Range("K25").Value = 15.01498053123629
$ws.Range("L25").Value = 9.267739841265991
$ws.Range("M25").Value = 17.2041700596169
$ws.Range("O25").Value = 22.11913391131329
